# Reformat the "optimization_parameters" sheet to match the current beta
# format: rename "Model" -> "production_function", add a new "L_curve"
# parameter row, drop the stray duplicated "value" header cells (C1:F1),
# and drop the obsolete "Deletion" row. Also re-point the active
# sheet/selection at this sheet (it becomes the active tab).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Row 1 header had "value" duplicated across C1:F1 (leftover/stray cells) -
# clear them so only A1/B1 remain.
$ws.Range("C1:F1").ClearContents()

# Insert a new row for the "L_curve" parameter right after "production_function"
# (old row 8 "Model"/"Sigmoid"), i.e. before old row 9 "estimate_params".
$ws.Rows.Item(9).Insert()

# Rename the "Model" label to "production_function" (row 8, same row/col).
$ws.Range("A8").Value = "production_function"

# Populate the newly inserted row 9 with the "L_curve" parameter.
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 0
$ws.Range("B9").NumberFormat = "0.00E+00"

# Remove the obsolete "Deletion" row (old row 16, shifted to row 17 after
# the insert above).
$ws.Rows.Item(17).Delete()

# This sheet is now the active tab; reflect that in the view state.
$ws.Activate()
$ws.Range("C1:F2").Select()
